$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value updates: Price (D) and/or Volume(1h) (E) columns ---
# D column values are forced to Text format first so that numeric-looking
# strings (e.g. "487.48") are not silently converted to real numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.470.50"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.934.97"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "487.48"
$ws.Range("E5").Value = "  +3.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.51"
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.737"
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000354"
$ws.Range("E11").Value = "  +3.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.11"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.73"
$ws.Range("E13").Value = "  +2.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.569.76"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.934.72"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.03"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.591.51"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "444.53"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.50"
$ws.Range("E22").Value = "  +4.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.89"
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.74"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.66"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.02"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.85"
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "719.43"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.73"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0920"
$ws.Range("E34").Value = "  +14.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "61.07"
$ws.Range("E37").Value = "  +5.69%  "
$ws.Range("E38").Value = "  -3.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.400"
$ws.Range("E39").Value = "  +18.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.97"
$ws.Range("E41").Value = "  +15.18%  "
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.16"
$ws.Range("E43").Value = "  +3.27%  "
$ws.Range("E44").Value = "  +6.51%  "
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "146.08"
$ws.Range("E50").Value = "  +0.67%  "

# --- Full row swaps: Coin (B), Link (C), Price (D), Volume(1h) (E) ---
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.23"
$ws.Range("E25").Value = "  +16.53%  "
$ws.Range("B26").Value = "Filecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.51"
$ws.Range("E26").Value = "  +13.40%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.25"
$ws.Range("E35").Value = "  +15.81%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "42.16"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.17"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0343"
$ws.Range("E51").Value = "  +30.88%  "
